$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal-looking string (e.g. "352.61")
# must be forced to Text, otherwise Excel/COM auto-converts them to numbers
# (dropping trailing zeros / precision, e.g. "70.30" -> 70.3).
# NumberFormat "@" forces text entry; Style="Normal" immediately resets the
# display format back to the original default so no visible formatting changes.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "51.544.62"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.791.22"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "352.61"
$ws.Range("E5").Value = "  -1.97%  "
Set-TextValue $ws.Range("D6") "108.72"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +5.22%  "
Set-TextValue $ws.Range("D10") "39.59"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("E12").Value = "  -1.93%  "
Set-TextValue $ws.Range("D13") "19.94"
$ws.Range("E13").Value = "  +2.10%  "
Set-TextValue $ws.Range("D14") "7.78"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D16").Value = "2.779.21"
Set-TextValue $ws.Range("D17") "0.932"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "51.524.47"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("E20").Value = "  +0.80%  "
Set-TextValue $ws.Range("D21") "13.31"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  -0.99%  "
Set-TextValue $ws.Range("D23") "70.30"
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue $ws.Range("D24") "266.57"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  -0.02%  "
Set-TextValue $ws.Range("D27") "25.90"
$ws.Range("E27").Value = "  -2.52%  "
Set-TextValue $ws.Range("D28") "0.165"
$ws.Range("E28").Value = "  +1.72%  "
Set-TextValue $ws.Range("D29") "10.29"
$ws.Range("E29").Value = "  -0.37%  "
Set-TextValue $ws.Range("D30") "37.02"
$ws.Range("E30").Value = "  +7.28%  "
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  +8.72%  "
Set-TextValue $ws.Range("D33") "52.22"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +8.18%  "
Set-TextValue $ws.Range("D35") "0.0443"
$ws.Range("E35").Value = "  -6.57%  "
Set-TextValue $ws.Range("D36") "0.0849"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +0.06%  "
Set-TextValue $ws.Range("D38") "18.52"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("E42").Value = "  -6.32%  "
Set-TextValue $ws.Range("D43") "119.97"
$ws.Range("E43").Value = "  +0.20%  "
Set-TextValue $ws.Range("D44") "21.85"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -2.94%  "
$ws.Range("D46").Value = "2.129.89"
$ws.Range("E46").Value = "  +1.96%  "
Set-TextValue $ws.Range("D47") "3.37"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("E48").Value = "  +6.07%  "
Set-TextValue $ws.Range("D49") "0.225"
$ws.Range("E49").Value = "  +17.70%  "
Set-TextValue $ws.Range("D50") "0.911"
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("E51").Value = "  +9.44%  "
